$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 updates
$ws.Range("H5").Value = 0.4
$ws.Range("N5").Value = 0.4767857142857143

# Row 6 updates
$ws.Range("E6").Value = 0.4285714285714285
$ws.Range("N6").Value = 0.353452380952381

# Row 8 updates
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0.2536904761904762
